$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.826.01'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  +1.98%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.804.93'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  +2.47%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '353.22'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '112.21'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +5.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.557'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  +2.72%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.626'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +8.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.25'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  +4.19%  '
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0840'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +1.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.93'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  +3.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.80'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +5.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.244.87'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  +2.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.803.63'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  +2.76%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.947'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  +4.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.812.34'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +2.10%  '
$ws.Range('E19').Value = '  +0.91%  '
$ws.Range('E20').Value = '  +7.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.58'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  +5.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0975'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +2.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.35'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +2.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.33'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +2.25%  '
$ws.Range('E25').Value = '  +2.68%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  +2.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.161'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +1.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '39.09'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +14.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.40'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +4.43%  '
$ws.Range('E31').Value = '  +1.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '52.18'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +1.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.11'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  +3.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0454'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +4.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.53'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +7.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.01'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  +3.77%  '
$ws.Range('E40').Value = '  +4.50%  '
$ws.Range('E41').Value = '  +2.72%  '
$ws.Range('E42').Value = '  +3.04%  '
$ws.Range('E43').Value = '  +1.75%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '119.68'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.96'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +1.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.53'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +11.18%  '
$ws.Range('E47').Value = '  +8.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.108.53'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +1.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.977'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  +8.01%  '
$ws.Range('E50').Value = '  +2.13%  '
$ws.Range('E51').Value = '  +7.92%  '
